$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row entry: row number, then new values for columns B,C,D,E,F,I,J,K,L,M,N (G and H are unchanged)
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$newData = @(
    @(2,1.02,1.073995936115403,1.070272582048248,1.076876138636622,1.084835956931483,1.059198335520961,1.078907638513878,1.072972861239862,1.079558917116555,1.087497946080979,1.080439811113003),
    @(3,1.02,1.075387383630167,1.071169776382003,1.078076691739842,1.086012962226143,1.059637956305715,1.079955918929432,1.073686185543798,1.080576121590518,1.088493168866402,1.081489580206987),
    @(4,1.02,1.076287080492588,1.071749751949517,1.078853072123655,1.086774237834576,1.059920766994024,1.080633065001921,1.074146548587213,1.081233286973565,1.089136235643914,1.082167687904569),
    @(5,1.02,1.07666515860288,1.071993438808648,1.079179355655758,1.087094202586794,1.060039266045107,1.08091746205031,1.074339798121822,1.081509314385537,1.089406365595682,1.082452488829387),
    @(6,1.02,1.076728630646252,1.072034347035971,1.079234133916215,1.087147921718549,1.060059139456266,1.08096519751455,1.074372228777632,1.081555646316374,1.089451709042615,1.082500292083462),
    @(7,1.02,1.076292132992449,1.071753008633482,1.078857432356267,1.086778513512967,1.059922351932949,1.080636866206799,1.07414913192315,1.08123697622446,1.089139845976793,1.082171494507594),
    @(8,1.02,1.074466322073606,1.070575911760885,1.077281967650727,1.085233800240257,1.059347250518186,1.07926215201903,1.07321418303896,1.079902901786248,1.087834475848692,1.080794828067995),
    @(9,1.02,1.071243745665187,1.068497294744742,1.0745021656434,1.082509216223612,1.058321129145695,1.076830704674656,1.071557383224358,1.077544056256524,1.085527183397322,1.078359927789502),
    @(10,1.02,1.069091556501673,1.067108486582854,1.072646343967482,1.080690923017359,1.057628416291614,1.075203488046156,1.070446500042143,1.075965923790539,1.083984094112197,1.076730400326746),
    @(11,1.02,1.068158672746015,1.066506369782714,1.071842090794726,1.079903095274195,1.057326398313075,1.074497360990493,1.069963947050343,1.075281216646645,1.083314725557782,1.07602327048967),
    @(12,1.02,1.067812006983975,1.066282601661097,1.07154325110588,1.079610383513947,1.057213902939458,1.074234840207082,1.069784473034654,1.075026677474047,1.083065908945479,1.075760376896635),
    @(13,1.02,1.06788637489282,1.066330605923406,1.071607358008664,1.079673174735627,1.057238047744431,1.074291162443132,1.069822981366714,1.075081286477038,1.083119289296682,1.07581677911672),
    @(14,1.02,1.068130020360654,1.066487875399674,1.071817390760221,1.079878901227026,1.057317105793799,1.074475665730402,1.069949116423652,1.075260180630869,1.08329416205584,1.07600154441982),
    @(15,1.02,1.068280118182911,1.066584759001336,1.071946784949045,1.080005645750532,1.057365774615422,1.074589313210496,1.070026801546821,1.075370375501679,1.083401882553577,1.076115353292372),
    @(16,1.02,1.069153448007667,1.067148431050574,1.072699705178386,1.08074319785932,1.057648416504022,1.075250318793814,1.070478493021144,1.076011336498706,1.084028492350923,1.076777297579436),
    @(17,1.02,1.069701000917812,1.06750180451197,1.073171810215099,1.081205710143373,1.057825155274118,1.075664537521849,1.070761415250509,1.076413026508975,1.084421224590208,1.077192104545577),
    @(18,1.02,1.070020285871359,1.067707848682052,1.073447116679814,1.08147543820674,1.057928044516949,1.075905996777711,1.070926291200732,1.076647193983297,1.084650183009578,1.077433906701278),
    @(19,1.02,1.070129138105207,1.067778092163447,1.073540978251016,1.08156740053835,1.057963093302025,1.075988303232257,1.070982484629527,1.076727016833662,1.08472823238546,1.077516330040431),
    @(20,1.02,1.069642263372155,1.067463898406953,1.073121164502541,1.081156091892179,1.057806213523591,1.075620111058868,1.070731075675712,1.076369942600578,1.084379100122394,1.077147614991923),
    @(21,1.02,1.068058277058931,1.06644156665229,1.071755544256414,1.079818322066298,1.057293833807095,1.074421340593186,1.069911979223274,1.07520750654245,1.083242671471236,1.075947142134679),
    @(22,1.02,1.067061484606933,1.065798119029336,1.070896319661534,1.078976763372216,1.056969871789644,1.073666271765002,1.069395635478089,1.074475429156913,1.082527091533012,1.075191001022151),
    @(23,1.02,1.067589988059565,1.066139286656083,1.07135186941374,1.079422933232038,1.05714178216515,1.074066677631266,1.069669487254458,1.074863632818307,1.082906535659054,1.075591975510672),
    @(24,1.02,1.069668804629009,1.067481026768463,1.073144049318938,1.081178512390875,1.057814773101654,1.075640185920824,1.070744785280247,1.076389410770479,1.084398134713426,1.077167718362485),
    @(25,1.02,1.072077510029953,1.069035200354865,1.075221259816905,1.083213910693068,1.058587921647529,1.077460380496658,1.071986817777866,1.078154843699845,1.086124525769187,1.078990497823396)
)

foreach ($entry in $newData) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $rowNum
        $ws.Range($cellRef).Value = $entry[$i + 1]
    }
}
